$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: was "OC" data, becomes "Sulfur" data ---
$ws.Range("D5").Value = "prefix|Emissions|Sulfur|Harmonized-DB"
$ws.Range("E5").Value = "Mt SO2/yr"
$ws.Range("F5").Value = 15
$ws.Range("G5").Value = 17.1428571428571

# --- Rows 6-8: was "Sulfur" data, becomes "OC" data ---
$ws.Range("D6").Value = "prefix|Emissions|OC|Harmonized-DB"
$ws.Range("E6").Value = "Mt OC/yr"
$ws.Range("F6").Formula = "=SUM(F7:F8)"

$ws.Range("D7").Value = "prefix|Emissions|OC|Harmonized-DB"
$ws.Range("E7").Value = "Mt OC/yr"
$ws.Range("F7").Value = 15

$ws.Range("D8").Value = "prefix|Emissions|OC|Harmonized-DB"
$ws.Range("E8").Value = "Mt OC/yr"
$ws.Range("F8").Value = 20

# --- Formatting: new cell style applied to the edited Unit/value columns ---
$ws.Range("E5:G5").Style = "Normal"
$ws.Range("E6:F6").Style = "Normal"
$ws.Range("E7:F7").Style = "Normal"
$ws.Range("E8:F8").Style = "Normal"

# --- Row heights ---
$ws.Rows("5").RowHeight = 13.8
$ws.Rows("6:8").RowHeight = 15.8
$ws.Rows("11:18").RowHeight = 15.8

# --- Extend used range down to row 11 without adding visible content ---
$ws.Range("G11").NumberFormat = "General"

# --- Selection follows the last-edited cell ---
[void]$ws.Range("G6").Select()
